# Apply "added hr_periods again" fix: rows 24-27 in columns B-D, F-H, J-L
# were re-aligned - a new "0.0" category row was restored at the top of
# each CKD-Stage group (rows 24-27), shifting the rest of the category
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, F and J hold category labels that look numeric ("0.0", "1.0",
# ...). Excel's COM layer auto-converts numeric-looking strings assigned
# via .Value into real numbers, which would lose the trailing ".0" text
# formatting. Force those ranges to Text format first so the values are
# stored verbatim as strings.
$ws.Range("B24:B27").NumberFormat = "@"
$ws.Range("F24:F26").NumberFormat = "@"
$ws.Range("J24:J27").NumberFormat = "@"

# Row 24
$ws.Range("B24").Value = "0.0"
$ws.Range("C24").Value = "1362 (90.6)"
$ws.Range("D24").Value = "6634 (93.8)"
$ws.Range("F24").Value = "0.0"
$ws.Range("G24").Value = "469 (90.0)"
$ws.Range("H24").Value = "2411 (93.0)"
$ws.Range("J24").Value = "0.0"
$ws.Range("K24").Value = "583 (91.5)"
$ws.Range("L24").Value = "2400 (93.1)"

# Row 25
$ws.Range("B25").Value = "1.0"
$ws.Range("C25").Value = "3 (0.2)"
$ws.Range("D25").Value = "1 (0.0)"
$ws.Range("F25").Value = "2.0"
$ws.Range("G25").Value = "7 (1.3)"
$ws.Range("H25").Value = "21 (0.8)"
$ws.Range("J25").Value = "1.0"
$ws.Range("K25").Value = "2 (0.3)"
$ws.Range("L25").Value = "1 (0.0)"

# Row 26
$ws.Range("B26").Value = "2.0"
$ws.Range("C26").Value = "14 (0.9)"
$ws.Range("D26").Value = "45 (0.6)"
$ws.Range("F26").Value = "3.0"
$ws.Range("G26").Value = "45 (8.6)"
$ws.Range("H26").Value = "161 (6.2)"
$ws.Range("J26").Value = "2.0"
$ws.Range("K26").Value = "8 (1.3)"
$ws.Range("L26").Value = "25 (1.0)"

# Row 27 (F/G/H belong to the "Diabetes Type" group here and are untouched)
$ws.Range("B27").Value = "3.0"
$ws.Range("C27").Value = "124 (8.3)"
$ws.Range("D27").Value = "393 (5.6)"
$ws.Range("J27").Value = "3.0"
$ws.Range("K27").Value = "44 (6.9)"
$ws.Range("L27").Value = "153 (5.9)"
